# repull data, push all data, mean calculation
# Updates column F ("dSF") values on Sheet1 for the rows whose
# dSF value changed after a data repull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -5
    4  = 1
    5  = 3
    6  = -3
    7  = -4
    8  = 1
    10 = 5
    11 = -2
    13 = 9
    15 = -5
    16 = -6
    19 = -1
    20 = -4
    21 = 4
    22 = 2
    23 = 3
    25 = -2
    26 = -2
    28 = -4
    29 = 1
    30 = -2
    32 = -1
    33 = -2
    36 = 1
    37 = 1
    38 = -3
    39 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
